$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Collapse the "- XXXX – Tarea XXXX" run sequence in the header row into a
#    single {{TEST_ID}} placeholder run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "- XXXX – Tarea XXXX"
$rng.Find.Execute() | Out-Null
$rng.Text = "{{TEST_ID}}"

# ---------------------------------------------------------------------------
# 2) Give the "Registro" row an explicit height (231 twips = 11.55 pt).
# ---------------------------------------------------------------------------
$table = $d.Tables(1)
$registroRow = $table.Rows(2)
$registroRow.Height = 11.55

# ---------------------------------------------------------------------------
# 3) Split the "Registro: " run into separate ":" and " " runs, and swap the
#    placeholder text for {{STEP_DESC}}.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Text = ": "
$rng2.Find.Execute() | Out-Null
$rng2.InsertAfter(" ")
$rng2.Collapse(1)
$rng2.MoveEnd(1, 1)
$rng2.Text = ":"

$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$rng3.Find.Text = "XXXXXXXXXXXXXXXXXXXXXXXXXX"
$rng3.Find.Execute() | Out-Null
$rng3.Text = "{{STEP_DESC}}"

# ---------------------------------------------------------------------------
# 4) Remove the empty "Ttulo1" spacer row that used to follow the Registro
#    row.
# ---------------------------------------------------------------------------
$table.Rows(3).Delete()

# ---------------------------------------------------------------------------
# 5) Move the lone _GoBack bookmark from after "Resumen del Problema:" to
#    right after the new {{STEP_DESC}} run.
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.ClearFormatting()
$rng4.Find.Text = "{{STEP_DESC}}"
$rng4.Find.Execute() | Out-Null
$afterStepDesc = $d.Range($rng4.End, $rng4.End)
$afterStepDesc.Bookmarks.Add("_GoBack") | Out-Null

# ---------------------------------------------------------------------------
# 6) Mark the "Default Paragraph Font" character style as semi-hidden.
# ---------------------------------------------------------------------------
$style = $d.Styles("Default Paragraph Font")
$style.Visibility = $false
